$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Wheel for Micro Continuous Rotation FS90R Servo -> now sourced from 3D Hubs ---
$ws.Range("B4").Value = "3D Hubs"
$ws.Range("D4").Value = 0.3
$ws.Range("E4").Value = 40
$ws.Range("F4").Formula = "=D4*E4"
$ws.Range("G4").Value = "https://www.3dhubs.com/"
$ws.Range("H4").Value = "Rev02 prints without support"

# --- Row 5: O-ring -> now a specific O-ring part number from Acklands Grainger ---
$ws.Range("A5").Value = "O-ring BS1806-022"
$ws.Range("D5").Value = 0.15
$ws.Range("E5").Value = 40
$ws.Range("F5").Formula = "=D5*E5"
$ws.Range("G5").Value = "https://www.acklandsgrainger.com/en/product/O-RING-70BUNA-1X1-16/p/HBS568-022"

# Copy the existing currency number formatting (used lower in the sheet) onto the
# newly-populated cost cells so they render as "US$ Unit Cost" / "US$ Subtotal" like the rest.
$ws.Range("D6").Copy()
$ws.Range("D4:D5").PasteSpecial(-4122)
$ws.Range("F6").Copy()
$ws.Range("F4:F5").PasteSpecial(-4122)

# --- Row 6: PCB from LPKF - add an estimate comment ---
$ws.Range("H6").Value = "Estimate"

# --- Row 7: Acrylic -> specific sheet size / updated pricing ---
$ws.Range("A7").Value = "Clear Acrylic 12""x24"""
$ws.Range("D7").Value = 7.18
$ws.Range("E7").Value = 1.5
$ws.Range("F7").Formula = "=D7*E7"

# --- Row 8: Acrylic cutting -> source added ---
$ws.Range("B8").Value = "Pololu"

# --- Header / title formatting ---
$ws.Range("A1").Font.Bold = $true
$ws.Range("A3:H3").Font.Bold = $true

$excel.Calculate()

# Restore the active selection like the authored workbook
$ws.Range("D9").Select()
